$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 odds updates
$ws.Range("G9").Value = 1.95
$ws.Range("H9").Value = 3.2
$ws.Range("I9").Value = 4.33
$ws.Range("J9").Value = 2.63
$ws.Range("K9").Value = 2.05
$ws.Range("L9").Value = 4.75
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 2.15
$ws.Range("R9").Value = 1.67
$ws.Range("U9").Value = 1.95
$ws.Range("V9").Value = 1.8
$ws.Range("X9").Value = 8.5
$ws.Range("Z9").Value = 17
$ws.Range("AA9").Value = 17
$ws.Range("AC9").Value = 8
$ws.Range("AF9").Value = 51
$ws.Range("AG9").Value = 351
$ws.Range("AH9").Value = 10
$ws.Range("AI9").Value = 21
$ws.Range("AJ9").Value = 15
$ws.Range("AL9").Value = 41
$ws.Range("AN9").Value = 3.75
$ws.Range("AO9").Value = 11
$ws.Range("AP9").Value = 23
$ws.Range("AU9").Value = 8.5
$ws.Range("AX9").Value = 6
$ws.Range("AY9").Value = 23

# Row 10 odds updates
$ws.Range("G10").Value = 1.85
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 3.7
$ws.Range("J10").Value = 2.5
$ws.Range("W10").Value = 8.5
$ws.Range("Z10").Value = 17
$ws.Range("AO10").Value = 10
$ws.Range("AQ10").Value = 34
$ws.Range("AS10").Value = 126
